$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2024-05-14 Tuesday"

# Update each table cell value (row-major order, 5 columns per row)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "41-24="
$t.Cell(1,2).Range.Text = "44+16="
$t.Cell(1,3).Range.Text = "11+9="
$t.Cell(1,4).Range.Text = "74-6="
$t.Cell(1,5).Range.Text = "5+33="
$t.Cell(2,1).Range.Text = "15+67="
$t.Cell(2,2).Range.Text = "16+60="
$t.Cell(2,3).Range.Text = "73-60="
$t.Cell(2,4).Range.Text = "26+62="
$t.Cell(2,5).Range.Text = "53-35="
$t.Cell(3,1).Range.Text = "62-1="
$t.Cell(3,2).Range.Text = "13+29="
$t.Cell(3,3).Range.Text = "32+12="
$t.Cell(3,4).Range.Text = "51-10="
$t.Cell(3,5).Range.Text = "34-24="
$t.Cell(4,1).Range.Text = "9+60="
$t.Cell(4,2).Range.Text = "75-32="
$t.Cell(4,3).Range.Text = "21-14="
$t.Cell(4,4).Range.Text = "4+15="
$t.Cell(4,5).Range.Text = "90-41="
$t.Cell(5,1).Range.Text = "65-47="
$t.Cell(5,2).Range.Text = "35+63="
$t.Cell(5,3).Range.Text = "72+25="
$t.Cell(5,4).Range.Text = "49-17="
$t.Cell(5,5).Range.Text = "87-44="
$t.Cell(6,1).Range.Text = "69+25="
$t.Cell(6,2).Range.Text = "42+32="
$t.Cell(6,3).Range.Text = "42+30="
$t.Cell(6,4).Range.Text = "10+34="
$t.Cell(6,5).Range.Text = "21+78="
$t.Cell(7,1).Range.Text = "91-77="
$t.Cell(7,2).Range.Text = "23-5="
$t.Cell(7,3).Range.Text = "73-3="
$t.Cell(7,4).Range.Text = "89-83="
$t.Cell(7,5).Range.Text = "36+10="
$t.Cell(8,1).Range.Text = "68-52="
$t.Cell(8,2).Range.Text = "63-60="
$t.Cell(8,3).Range.Text = "91-32="
$t.Cell(8,4).Range.Text = "33+20="
$t.Cell(8,5).Range.Text = "98-89="
$t.Cell(9,1).Range.Text = "3+45="
$t.Cell(9,2).Range.Text = "92-25="
$t.Cell(9,3).Range.Text = "87-23="
$t.Cell(9,4).Range.Text = "2+61="
$t.Cell(9,5).Range.Text = "71+26="
$t.Cell(10,1).Range.Text = "3+52="
$t.Cell(10,2).Range.Text = "56+43="
$t.Cell(10,3).Range.Text = "99-92="
$t.Cell(10,4).Range.Text = "51-30="
$t.Cell(10,5).Range.Text = "68-33="
$t.Cell(11,1).Range.Text = "49-8="
$t.Cell(11,2).Range.Text = "10+18="
$t.Cell(11,3).Range.Text = "14-12="
$t.Cell(11,4).Range.Text = "24-19="
$t.Cell(11,5).Range.Text = "27+3="
$t.Cell(12,1).Range.Text = "68+7="
$t.Cell(12,2).Range.Text = "55-13="
$t.Cell(12,3).Range.Text = "5+28="
$t.Cell(12,4).Range.Text = "38-28="
$t.Cell(12,5).Range.Text = "72-11="
$t.Cell(13,1).Range.Text = "35-12="
$t.Cell(13,2).Range.Text = "52+32="
$t.Cell(13,3).Range.Text = "78-32="
$t.Cell(13,4).Range.Text = "47+52="
$t.Cell(13,5).Range.Text = "45-25="
$t.Cell(14,1).Range.Text = "74+4="
$t.Cell(14,2).Range.Text = "26+42="
$t.Cell(14,3).Range.Text = "2+40="
$t.Cell(14,4).Range.Text = "93-53="
$t.Cell(14,5).Range.Text = "65+16="
$t.Cell(15,1).Range.Text = "5+34="
$t.Cell(15,2).Range.Text = "80+3="
$t.Cell(15,3).Range.Text = "57-44="
$t.Cell(15,4).Range.Text = "6+57="
$t.Cell(15,5).Range.Text = "38+61="
$t.Cell(16,1).Range.Text = "71+17="
$t.Cell(16,2).Range.Text = "0+59="
$t.Cell(16,3).Range.Text = "99-13="
$t.Cell(16,4).Range.Text = "26+39="
$t.Cell(16,5).Range.Text = "52+24="
$t.Cell(17,1).Range.Text = "1+78="
$t.Cell(17,2).Range.Text = "64-5="
$t.Cell(17,3).Range.Text = "12-8="
$t.Cell(17,4).Range.Text = "61-33="
$t.Cell(17,5).Range.Text = "98-12="
$t.Cell(18,1).Range.Text = "12+38="
$t.Cell(18,2).Range.Text = "23-19="
$t.Cell(18,3).Range.Text = "43+31="
$t.Cell(18,4).Range.Text = "38-30="
$t.Cell(18,5).Range.Text = "44+21="
$t.Cell(19,1).Range.Text = "11+87="
$t.Cell(19,2).Range.Text = "32-25="
$t.Cell(19,3).Range.Text = "76+20="
$t.Cell(19,4).Range.Text = "38+57="
$t.Cell(19,5).Range.Text = "91-90="
$t.Cell(20,1).Range.Text = "34-12="
$t.Cell(20,2).Range.Text = "67+31="
$t.Cell(20,3).Range.Text = "33+30="
$t.Cell(20,4).Range.Text = "78-71="
$t.Cell(20,5).Range.Text = "95-10="
